$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -78.298
$ws.Range("B2").Value = -78.1288

$ws.Range("A3").Value = 33.6418
$ws.Range("B3").Value = 33.7808

$ws.Range("A4").Value = -75.4785
$ws.Range("B4").Value = -75.6532

$ws.Range("A5").Value = 36.6811
$ws.Range("B5").Value = 36.5425
